$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. D2 (row 7) LED description: "Red ..." -> the correct "Orange ..." text.
#    (The BOM previously had the wrong/duplicated "Red" description copied
#    into the ORANGE LED row; fix it to the proper Orange description.)
# ---------------------------------------------------------------------------
$ws.Range("D7").Value = "Orange 0603 130` Clear 54 mcd 2 V Surface Mount ChipLED ;"

# ---------------------------------------------------------------------------
# 2. Merge R8 (2k2, row 18) into R10 (1k, row 13): R8 was fixed from 2k2 to
#    1k, so it is now the same part as R10 and both reference designators
#    are combined onto a single BOM line ("R8,R10"). Capture the R10 row,
#    remove it from its old (alphabetically sorted) position, and re-insert
#    the merged row where "R8,R10" sorts alphabetically (after R7, before U1).
# ---------------------------------------------------------------------------
$srcRow = 13
$rowVals = @()
for ($c = 1; $c -le 7; $c++) {
    $rowVals += $ws.Cells.Item($srcRow, $c).Value2
}

$ws.Rows($srcRow).Delete()

$destRow = 17
$ws.Rows($destRow).Insert()

for ($c = 1; $c -le 7; $c++) {
    $ws.Cells.Item($destRow, $c).Value = $rowVals[$c - 1]
}
$ws.Cells.Item($destRow, 3).Value = "R8,R10"

# Remove the now-redundant old R8 (2k2) line entirely.
$ws.Rows(18).Delete()

# ---------------------------------------------------------------------------
# 3. Add hyperlinks on the OCTOPART_URL column for the rows Excel highlighted
#    (D1 / D2 / the merged R8,R10 row), turning the plain URL text into a
#    real hyperlink (keeping the same displayed text).
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("G7"), $ws.Range("G7").Text)
$ws.Hyperlinks.Add($ws.Range("G6"), $ws.Range("G6").Text)
$ws.Hyperlinks.Add($ws.Range("G17"), $ws.Range("G17").Text)

# ---------------------------------------------------------------------------
# 4. The defined name / print range shrank by one row (22 -> 21 rows).
# ---------------------------------------------------------------------------
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!H0FR70") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$H`$21"
    }
}

# ---------------------------------------------------------------------------
# 5. Update the sort-state helper range/condition to match the new extent.
# ---------------------------------------------------------------------------
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("C4"))
$ws.Sort.SetRange($ws.Range("A2:G22"))
$ws.Sort.Apply()

# ---------------------------------------------------------------------------
# 6. Restore the view/selection state left behind after the edit.
# ---------------------------------------------------------------------------
$ws.Range("G17").Select()

Write-Host "edit complete"
